$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 = Disease Ontology ("do") -> Disease Ontology release date
$ws.Range("E3").Value = "v2023-04-01"

# Row 4 = Experimental Factor Ontology ("efo") -> EFO version
$ws.Range("E4").Value = "v3.52.0"

# Move the active selection to E4
$ws.Range("E4").Select()
